# Add new "Teste" rows (41-59) to the statistics sheet and convert the
# previously-text numeric cells on row 40 to real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$H_STD = "HTHG,HTAG,HTR,HS,AS,HST,AST,HC,AC,HY,AY,HR,AR"

# --- Row 40: C/D/E/F were stored as text ("0.5","0.2","100","1000"); turn
# them into real numeric cells (same values). G40/H40 stay as they were.
$ws.Cells.Item(40, 3).Value = 0.5
$ws.Cells.Item(40, 4).Value = 0.2
$ws.Cells.Item(40, 5).Value = 100
$ws.Cells.Item(40, 6).Value = 1000

# --- Rows 41-58: new numeric "Teste" rows.
$newRows = @(
    @{ Row = 41; B = "E0"; C = 0.5; D = 0.2; E = 100; F = 10000; G = "8 de 10" },
    @{ Row = 42; B = "D1"; C = 0.8; D = 0.2; E = 100; F = 10000; G = "7 de 10" },
    @{ Row = 43; B = "D1"; C = 0.8; D = 0.2; E = 100; F = 10000; G = "4 de 10" },
    @{ Row = 44; B = "D1"; C = 0.8; D = 0.2; E = 100; F = 10000; G = "4 de 10" },
    @{ Row = 45; B = "D1"; C = 0.8; D = 0.2; E = 100; F = 10000; G = "6 de 10" },
    @{ Row = 46; B = "D1"; C = 0.8; D = 0.2; E = 100; F = 10000; G = "6 de 10" },
    @{ Row = 47; B = "D1"; C = 0.8; D = 0.2; E = 100; F = 10000; G = "6 de 10" },
    @{ Row = 48; B = "D1"; C = 0.8; D = 0.2; E = 100; F = 10000; G = "4 de 10" },
    @{ Row = 49; B = "D1"; C = 0.8; D = 0.2; E = 100; F = 10000; G = "5 de 10" },
    @{ Row = 50; B = "D1"; C = 0.8; D = 0.2; E = 100; F = 10000; G = "6 de 10" },
    @{ Row = 51; B = "D1"; C = 0.8; D = 0.2; E = 100; F = 10000; G = "7 de 10" },
    @{ Row = 52; B = "D1"; C = 0.8; D = 0.2; E = 100; F = 10000; G = "4 de 10" },
    @{ Row = 53; B = "D1"; C = 0.8; D = 0.2; E = 100; F = 10000; G = "5 de 10" },
    @{ Row = 54; B = "D1"; C = 0.8; D = 0.2; E = 100; F = 10000; G = "6 de 10" },
    @{ Row = 55; B = "D1"; C = 0.8; D = 0.2; E = 100; F = 10000; G = "5 de 10" },
    @{ Row = 56; B = "D1"; C = 0.8; D = 0.2; E = 100; F = 10000; G = "7 de 10" },
    @{ Row = 57; B = "D1"; C = 0.8; D = 0.2; E = 100; F = 10000; G = "4 de 10" },
    @{ Row = 58; B = "D1"; C = 0.8; D = 0.2; E = 100; F = 10000; G = "5 de 10" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = "Teste"
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $H_STD
}

# --- Row 59: new row whose C/D/E/F are stored as text (like row 40 used to
# be before this edit), using a leading apostrophe so Excel keeps them as
# text instead of coercing the numeric-looking strings into numbers.
$ws.Cells.Item(59, 1).Value = "Teste"
$ws.Cells.Item(59, 2).Value = "E0"
$ws.Cells.Item(59, 3).Value = "'0.8"
$ws.Cells.Item(59, 4).Value = "'0"
$ws.Cells.Item(59, 5).Value = "'100"
$ws.Cells.Item(59, 6).Value = "'100000"
$ws.Cells.Item(59, 7).Value = "8 de 10"
$ws.Cells.Item(59, 8).Value = $H_STD
